$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are forced to text (matching the source inlineStr cells) by
# temporarily applying a text number format, then resetting the style back to
# "Normal" so no explicit style index is left on the cell.
function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws 'D2' '64.090.52'
$ws.Range('E2').Value = '  +0.92%  '
Set-TextValue $ws 'D3' '3.067.52'
$ws.Range('E3').Value = '  +0.18%  '
Set-TextValue $ws 'D4' '1.00'
$ws.Range('E4').Value = '  -0.11%  '
Set-TextValue $ws 'D5' '558.26'
$ws.Range('E5').Value = '  +1.65%  '
Set-TextValue $ws 'D6' '146.16'
$ws.Range('E6').Value = '  +5.06%  '
Set-TextValue $ws 'D7' '1.00'
$ws.Range('E7').Value = '  +0.09%  '
Set-TextValue $ws 'D8' '3.064.98'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws 'D10' '0.153'
$ws.Range('E10').Value = '  +2.49%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws 'D11' '6.25'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('E12').Value = '  +3.84%  '
Set-TextValue $ws 'D13' '0.0000228'
$ws.Range('E13').Value = '  +0.54%  '
Set-TextValue $ws 'D14' '35.20'
$ws.Range('E14').Value = '  +1.66%  '
Set-TextValue $ws 'D15' '3.568.94'
$ws.Range('E15').Value = '  +0.16%  '
Set-TextValue $ws 'D16' '64.097.23'
$ws.Range('E16').Value = '  +0.80%  '
Set-TextValue $ws 'D17' '3.068.67'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('E18').Value = '  +1.15%  '
Set-TextValue $ws 'D19' '6.78'
$ws.Range('E19').Value = '  +0.98%  '
Set-TextValue $ws 'D20' '476.93'
Set-TextValue $ws 'D21' '13.93'
$ws.Range('E21').Value = '  +2.79%  '
Set-TextValue $ws 'D22' '0.675'
$ws.Range('E22').Value = '  -0.40%  '
Set-TextValue $ws 'D23' '7.54'
$ws.Range('E23').Value = '  +4.91%  '
Set-TextValue $ws 'D24' '13.51'
$ws.Range('E24').Value = '  +7.91%  '
Set-TextValue $ws 'D25' '81.57'
$ws.Range('E25').Value = '  +0.45%  '
Set-TextValue $ws 'D26' '0.999'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  +1.85%  '
Set-TextValue $ws 'D28' '8.12'
$ws.Range('E28').Value = '  +2.12%  '
Set-TextValue $ws 'D29' '2.06'
$ws.Range('E29').Value = '  +4.11%  '
Set-TextValue $ws 'D30' '1.00'
$ws.Range('E30').Value = '  +0.09%  '
Set-TextValue $ws 'D31' '26.21'
$ws.Range('E31').Value = '  +0.89%  '
Set-TextValue $ws 'D32' '1.16'
$ws.Range('E32').Value = '  +0.93%  '
Set-TextValue $ws 'D33' '2.48'
$ws.Range('E33').Value = '  +2.74%  '
Set-TextValue $ws 'D34' '5.58'
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('E35').Value = '  +3.54%  '
Set-TextValue $ws 'D36' '54.81'
$ws.Range('E36').Value = '  -1.76%  '
Set-TextValue $ws 'D37' '461.74'
$ws.Range('E37').Value = '  -1.17%  '
$ws.Range('E38').Value = '  +18.47%  '
Set-TextValue $ws 'D39' '0.0831'
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('E40').Value = '  +2.86%  '
Set-TextValue $ws 'D41' '2.965.15'
$ws.Range('E41').Value = '  -5.24%  '
Set-TextValue $ws 'D42' '8.27'
$ws.Range('E42').Value = '  +0.49%  '
$ws.Range('E43').Value = '  -3.83%  '
Set-TextValue $ws 'D44' '28.00'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('E45').Value = '  +4.21%  '
Set-TextValue $ws 'D46' '2.15'
$ws.Range('E46').Value = '  +5.00%  '
$ws.Range('E48').Value = '  +2.64%  '
Set-TextValue $ws 'D49' '119.86'
Set-TextValue $ws 'D50' '0.0₃0517'
$ws.Range('E50').Value = '  +1.41%  '
Set-TextValue $ws 'D51' '2.08'
